{"js": "// The document text is being rearranged: several paragraphs under\n// \"Objetivos\", \"Docente(s) Respons\u00e1vel(eis)\", \"Programa resumido\",\n// \"Programa\", \"Avalia\u00e7\u00e3o\" (M\u00e9todo/Crit\u00e9rio/Norma de recupera\u00e7\u00e3o) and\n// \"Bibliografia\" swap their textual content in a cyclic fashion while\n// every paragraph keeps its own style/formatting in place.\n//\n// Because the new text for one slot is the old text that used to live\n// in a *different* slot (a rotation), we address every target purely\n// by its paragraph/run position rather than by searching for its old\n// text - that way there is no risk of one replacement accidentally\n// matching text that a previous step just inserted.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Simple single-run paragraphs: replace the whole paragraph's text\n// while leaving its paragraph/run formatting untouched.\nconst simpleReplacements = [\n  { index: 5, text: \"A definir, de acordo com o t\u00f3pico programado.\" },\n  { index: 6, text: \"To be defined, according to the programmed topic.\" },\n  {\n    index: 8,\n    text:\n      \"Complementar a forma\u00e7\u00e3o dos estudantes abordando, com maior profundidade, t\u00f3picos atuais e relevantes e atualizar com temas no estado da arte.\",\n  },\n  {\n    index: 10,\n    text:\n      \"O conte\u00fado desta disciplina optativa ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares ao conte\u00fado regular do curso de gradua\u00e7\u00e3o.\",\n  },\n  {\n    index: 11,\n    text:\n      \"Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics.\",\n  },\n  {\n    index: 13,\n    text:\n      \"Este curso dever\u00e1 conter duas avalia\u00e7\u00f5es escritas denominadas P1 e P2. A P2 dever\u00e1 englobar toda a mat\u00e9ria ministrada ao longo do semestre, abrangendo todos os t\u00f3picos previstos na ementa.\",\n  },\n  {\n    index: 18,\n    text: \"5840730 - Antonio Jefferson da Silva Machado\",\n  },\n];\n\nfor (const { index, text } of simpleReplacements) {\n  const range = paragraphs.items[index].getRange();\n  range.insertText(text, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// The \"Avalia\u00e7\u00e3o\" paragraph (index 16) holds three bold labels\n// (\"M\u00e9todo: \", \"Crit\u00e9rio: \", \"Norma de recupera\u00e7\u00e3o: \") each followed by\n// its own run of plain text + a line break. Only the plain-text runs\n// change. Locate each one (in this still-unmodified paragraph) by its\n// current, unique text before making any edits, then replace all three\n// - this sidesteps the fact that the new values for these runs are\n// text that used to belong to other paragraphs (a rotation).\nconst avaliacaoRange = paragraphs.items[16].getRange();\n\nconst metodoValue = avaliacaoRange.search(\n  \"Este curso dever\u00e1 conter duas avalia\u00e7\u00f5es escritas denominadas P1 e P2. A P2 dever\u00e1 englobar toda a mat\u00e9ria ministrada ao longo do semestre, abrangendo todos os t\u00f3picos previstos na ementa.\",\n  { matchCase: true }\n);\nconst criterioValue = avaliacaoRange.search(\n  \"A m\u00e9dia do semestre ser\u00e1 computada com base na rela\u00e7\u00e3o: M=(P1+2P2)/3\",\n  { matchCase: true }\n);\nconst normaValue = avaliacaoRange.search(\n  \"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre. A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 computada com base na rela\u00e7\u00e3o: MF=(M+RC)/2\",\n  { matchCase: true }\n);\nmetodoValue.load(\"items\");\ncriterioValue.load(\"items\");\nnormaValue.load(\"items\");\nawait context.sync();\n\nmetodoValue.items[0].insertText(\n  \"A m\u00e9dia do semestre ser\u00e1 computada com base na rela\u00e7\u00e3o: M=(P1+2P2)/3\",\n  Word.InsertLocation.replace\n);\ncriterioValue.items[0].insertText(\n  \"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre. A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 computada com base na rela\u00e7\u00e3o: MF=(M+RC)/2\",\n  Word.InsertLocation.replace\n);\nnormaValue.items[0].insertText(\n  \"Apostila ou texto fornecido pelo docente respons\u00e1vel. Artigos extra\u00eddos de revistas especializadas nas \u00e1reas de Ci\u00eancias e Tecnologia.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# The document text is being rearranged: several paragraphs under\n# \"Objetivos\", \"Docente(s) Responsavel(eis)\", \"Programa resumido\",\n# \"Programa\", \"Avaliacao\" (Metodo/Criterio/Norma de recuperacao) and\n# \"Bibliografia\" swap their textual content in a cyclic fashion while\n# every paragraph keeps its own style/formatting in place.\n#\n# Because the new text for one slot is the old text that used to live\n# in a *different* slot (a rotation), every target is addressed purely\n# by its paragraph position (1-based, Word COM style) rather than by\n# searching for its old text - that way there is no risk of one\n# replacement accidentally matching text that a previous step just\n# inserted.\n\n$d = $word.ActiveDocument\n\n# Single-run paragraphs: replacing Range.Text in place keeps the\n# paragraph's own formatting (style, italics, ...) untouched.\n$d.Paragraphs(6).Range.Text = \"A definir, de acordo com o t\u00f3pico programado.\"\n$d.Paragraphs(7).Range.Text = \"To be defined, according to the programmed topic.\"\n$d.Paragraphs(9).Range.Text = \"Complementar a forma\u00e7\u00e3o dos estudantes abordando, com maior profundidade, t\u00f3picos atuais e relevantes e atualizar com temas no estado da arte.\"\n$d.Paragraphs(11).Range.Text = \"O conte\u00fado desta disciplina optativa ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares ao conte\u00fado regular do curso de gradua\u00e7\u00e3o.\"\n$d.Paragraphs(12).Range.Text = \"Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics.\"\n$d.Paragraphs(14).Range.Text = \"Este curso dever\u00e1 conter duas avalia\u00e7\u00f5es escritas denominadas P1 e P2. A P2 dever\u00e1 englobar toda a mat\u00e9ria ministrada ao longo do semestre, abrangendo todos os t\u00f3picos previstos na ementa.\"\n$d.Paragraphs(19).Range.Text = \"5840730 - Antonio Jefferson da Silva Machado\"\n\n# The \"Avalia\u00e7\u00e3o\" paragraph (#17) holds three bold labels (\"M\u00e9todo: \",\n# \"Crit\u00e9rio: \", \"Norma de recupera\u00e7\u00e3o: \"), each followed by its own run\n# of plain text + a line break. Only the plain-text runs change.\n# Locate all three (in the still-unmodified paragraph) first, then\n# apply the edits starting from the rightmost match so that replacing\n# one run's text never shifts the character offsets of a match found\n# earlier in the paragraph.\n$avaliacaoParagraph = $d.Paragraphs(17)\n\n$seek1 = $avaliacaoParagraph.Range.Duplicate\n$seek1.Find.Execute(\"Este curso dever\u00e1 conter duas avalia\u00e7\u00f5es escritas denominadas P1 e P2. A P2 dever\u00e1 englobar toda a mat\u00e9ria ministrada ao longo do semestre, abrangendo todos os t\u00f3picos previstos na ementa.\") | Out-Null\n$metodoValueRange = $seek1.Duplicate\n\n$seek2 = $avaliacaoParagraph.Range.Duplicate\n$seek2.Find.Execute(\"A m\u00e9dia do semestre ser\u00e1 computada com base na rela\u00e7\u00e3o: M=(P1+2P2)/3\") | Out-Null\n$criterioValueRange = $seek2.Duplicate\n\n$seek3 = $avaliacaoParagraph.Range.Duplicate\n$seek3.Find.Execute(\"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre. A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 computada com base na rela\u00e7\u00e3o: MF=(M+RC)/2\") | Out-Null\n$normaValueRange = $seek3.Duplicate\n\n$normaValueRange.Text = \"Apostila ou texto fornecido pelo docente respons\u00e1vel. Artigos extra\u00eddos de revistas especializadas nas \u00e1reas de Ci\u00eancias e Tecnologia.\"\n$criterioValueRange.Text = \"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre. A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 computada com base na rela\u00e7\u00e3o: MF=(M+RC)/2\"\n$metodoValueRange.Text = \"A m\u00e9dia do semestre ser\u00e1 computada com base na rela\u00e7\u00e3o: M=(P1+2P2)/3\"\n"}
